$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency snapshot: updated prices/volume percentages for
# existing rows and shifted rows 45-51 down to make room for a new coin
# (EnergySwap) inserted at row 45.
# Each cell value is prefixed with a quote-prefix character ("'") so that
# Excel stores it as literal text and does not auto-convert numeric-looking
# strings (e.g. "1.00", "536.59") into numbers, matching the source data which
# stores these values as text.
$updates = @(
    @{ Cell = 'D2'; Value = "'59.959.81" }
    @{ Cell = 'E2'; Value = "'  +1.43%  " }
    @{ Cell = 'D3'; Value = "'2.653.79" }
    @{ Cell = 'E3'; Value = "'  +2.40%  " }
    @{ Cell = 'E4'; Value = "'  +0.01%  " }
    @{ Cell = 'D5'; Value = "'536.59" }
    @{ Cell = 'E5'; Value = "'  +1.37%  " }
    @{ Cell = 'D6'; Value = "'145.91" }
    @{ Cell = 'E6'; Value = "'  +4.27%  " }
    @{ Cell = 'D7'; Value = "'0.998" }
    @{ Cell = 'E7'; Value = "'  +0.03%  " }
    @{ Cell = 'E8'; Value = "'  +1.52%  " }
    @{ Cell = 'D9'; Value = "'2.672.34" }
    @{ Cell = 'E9'; Value = "'  +2.65%  " }
    @{ Cell = 'D10'; Value = "'6.69" }
    @{ Cell = 'E10'; Value = "'  +3.93%  " }
    @{ Cell = 'E11'; Value = "'  +2.22%  " }
    @{ Cell = 'E12'; Value = "'  +1.68%  " }
    @{ Cell = 'E13'; Value = "'  -1.25%  " }
    @{ Cell = 'D14'; Value = "'3.129.83" }
    @{ Cell = 'E14'; Value = "'  +2.45%  " }
    @{ Cell = 'D15'; Value = "'59.884.09" }
    @{ Cell = 'E15'; Value = "'  +1.40%  " }
    @{ Cell = 'D16'; Value = "'21.27" }
    @{ Cell = 'E16'; Value = "'  +4.03%  " }
    @{ Cell = 'D17'; Value = "'2.675.93" }
    @{ Cell = 'E17'; Value = "'  +3.29%  " }
    @{ Cell = 'E18'; Value = "'  +1.63%  " }
    @{ Cell = 'D19'; Value = "'345.13" }
    @{ Cell = 'E19'; Value = "'  -0.64%  " }
    @{ Cell = 'D20'; Value = "'4.43" }
    @{ Cell = 'E20'; Value = "'  +2.29%  " }
    @{ Cell = 'D21'; Value = "'10.32" }
    @{ Cell = 'E21'; Value = "'  +2.29%  " }
    @{ Cell = 'D22'; Value = "'6.39" }
    @{ Cell = 'E22'; Value = "'  -0.29%  " }
    @{ Cell = 'E23'; Value = "'  +0.04%  " }
    @{ Cell = 'D24'; Value = "'67.49" }
    @{ Cell = 'E24'; Value = "'  +0.22%  " }
    @{ Cell = 'E25'; Value = "'  +2.61%  " }
    @{ Cell = 'E26'; Value = "'  -0.14%  " }
    @{ Cell = 'D27'; Value = "'1.00" }
    @{ Cell = 'E27'; Value = "'  +0.27%  " }
    @{ Cell = 'D28'; Value = "'7.36" }
    @{ Cell = 'E28'; Value = "'  +2.84%  " }
    @{ Cell = 'E29'; Value = "'  +2.79%  " }
    @{ Cell = 'D30'; Value = "'0.998" }
    @{ Cell = 'E30'; Value = "'  -0.06%  " }
    @{ Cell = 'E31'; Value = "'  +3.05%  " }
    @{ Cell = 'D32'; Value = "'5.90" }
    @{ Cell = 'E32'; Value = "'  +0.71%  " }
    @{ Cell = 'D33'; Value = "'19.12" }
    @{ Cell = 'E33'; Value = "'  +1.89%  " }
    @{ Cell = 'D34'; Value = "'150.32" }
    @{ Cell = 'E34'; Value = "'  +1.09%  " }
    @{ Cell = 'D35'; Value = "'4.06" }
    @{ Cell = 'E35'; Value = "'  +2.22%  " }
    @{ Cell = 'E36'; Value = "'  +3.77%  " }
    @{ Cell = 'E37'; Value = "'  +0.46%  " }
    @{ Cell = 'E38'; Value = "'  +2.00%  " }
    @{ Cell = 'E39'; Value = "'  -0.15%  " }
    @{ Cell = 'D40'; Value = "'294.10" }
    @{ Cell = 'E40'; Value = "'  +9.19%  " }
    @{ Cell = 'E41'; Value = "'  +2.31%  " }
    @{ Cell = 'E42'; Value = "'  +0.09%  " }
    @{ Cell = 'E43'; Value = "'  +1.94%  " }
    @{ Cell = 'D44'; Value = "'0.0546" }
    @{ Cell = 'E44'; Value = "'  +5.40%  " }
    @{ Cell = 'B45'; Value = "'EnergySwap" }
    @{ Cell = 'C45'; Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = 'D45'; Value = "'19.41" }
    @{ Cell = 'E45'; Value = "'  +5.49%  " }
    @{ Cell = 'B46'; Value = "'WhiteBITCoin" }
    @{ Cell = 'C46'; Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt" }
    @{ Cell = 'D46'; Value = "'10.74" }
    @{ Cell = 'E46'; Value = "'  +0.00%  " }
    @{ Cell = 'B47'; Value = "'Stellar" }
    @{ Cell = 'C47'; Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" }
    @{ Cell = 'D47'; Value = "'0.0958" }
    @{ Cell = 'E47'; Value = "'  -0.08%  " }
    @{ Cell = 'D48'; Value = "'0.0228" }
    @{ Cell = 'E48'; Value = "'  +2.86%  " }
    @{ Cell = 'B49'; Value = "'Maker" }
    @{ Cell = 'C49'; Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" }
    @{ Cell = 'D49'; Value = "'1.976.08" }
    @{ Cell = 'E49'; Value = "'  +1.23%  " }
    @{ Cell = 'B50'; Value = "'InjectiveProtocol" }
    @{ Cell = 'C50'; Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj" }
    @{ Cell = 'D50'; Value = "'18.57" }
    @{ Cell = 'E50'; Value = "'  +1.95%  " }
    @{ Cell = 'B51'; Value = "'RenderToken" }
    @{ Cell = 'C51'; Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" }
    @{ Cell = 'D51'; Value = "'4.58" }
    @{ Cell = 'E51'; Value = "'  +0.71%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
